$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 6 (item 4564)
$ws.Range("H6").Value = 3277.75
$ws.Range("I6").Value = 3277.75
$ws.Range("K6").Value = 9833.25
$ws.Range("M6").Value = -9721.25

# row 28 (item 27772)
$ws.Range("H28").Value = 463186.72
$ws.Range("I28").Value = 505199.1
$ws.Range("K28").Value = 505199.1
$ws.Range("M28").Value = -504714.1

# row 53 (item 5479)
$ws.Range("H53").Value = 341.85715
$ws.Range("I53").Value = 230.75
$ws.Range("J53").Value = 490
$ws.Range("K53").Value = 230.75
$ws.Range("L53").Value = 490
$ws.Range("M53").Value = 406.25
$ws.Range("N53").Value = -1764

# row 62 (item 27781)
$ws.Range("H62").Value = 3595421.8
$ws.Range("I62").Value = 5857672.5
$ws.Range("J62").Value = 13525
$ws.Range("K62").Value = 5857672.5
$ws.Range("L62").Value = 13525
$ws.Range("M62").Value = -5857048.5
$ws.Range("N62").Value = -14773

# row 65 (item 27781)
$ws.Range("H65").Value = 3595421.8
$ws.Range("I65").Value = 5857672.5
$ws.Range("J65").Value = 13525
$ws.Range("K65").Value = 29288362.5
$ws.Range("L65").Value = 67625
$ws.Range("M65").Value = -29285242.5
$ws.Range("N65").Value = -73865

# row 76 (item 12602)
$ws.Range("H76").Value = 3475347
$ws.Range("I76").Value = 4118189.2
$ws.Range("J76").Value = 3998
$ws.Range("K76").Value = 4118189.2
$ws.Range("L76").Value = 3998
$ws.Range("M76").Value = -4117874.2
$ws.Range("N76").Value = -4628

# row 79 (item 12602)
$ws.Range("H79").Value = 3475347
$ws.Range("I79").Value = 4118189.2
$ws.Range("J79").Value = 3998
$ws.Range("K79").Value = 4118189.2
$ws.Range("L79").Value = 3998
$ws.Range("M79").Value = -4117097.2
$ws.Range("N79").Value = -6182

# row 98 (item 36237)
$ws.Range("H98").Value = 400547
$ws.Range("I98").Value = 466728.9
$ws.Range("J98").Value = 3455.5
$ws.Range("K98").Value = 466728.9
$ws.Range("L98").Value = 3455.5
$ws.Range("M98").Value = -465230.9
$ws.Range("N98").Value = -6451.5

# row 111 (item 27768)
$ws.Range("H111").Value = 1030.5333
$ws.Range("I111").Value = 539.8333
$ws.Range("J111").Value = 2993.3333
$ws.Range("K111").Value = 1619.4999
$ws.Range("L111").Value = 8979.999899999999
$ws.Range("M111").Value = 1447.5001
$ws.Range("N111").Value = -15113.9999

# row 122 (item 36237)
$ws.Range("H122").Value = 400547
$ws.Range("I122").Value = 466728.9
$ws.Range("J122").Value = 3455.5
$ws.Range("K122").Value = 1400186.7
$ws.Range("L122").Value = 10366.5
$ws.Range("M122").Value = -1397736.7
$ws.Range("N122").Value = -15266.5

# row 138 (item 44169)
$ws.Range("H138").Value = 2147.1
$ws.Range("I138").Value = 2006.762
$ws.Range("K138").Value = 6020.286
$ws.Range("M138").Value = -880.2860000000001

$ws = $wb.Worksheets.Item("ARM")
# row 6 (item 2226)
$ws.Range("H6").Value = 38001.2
$ws.Range("I6").Value = 38001.2
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 38001.2
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -37828.2
$ws.Range("N6").ClearContents()

# row 32 (item 44147)
$ws.Range("H32").Value = 14708.619
$ws.Range("I32").Value = 2742.0667
$ws.Range("J32").Value = 114429.89
$ws.Range("K32").Value = 2742.0667
$ws.Range("L32").Value = 114429.89
$ws.Range("M32").Value = -2455.0667
$ws.Range("N32").Value = -115003.89

# row 63 (item 12528)
$ws.Range("H63").Value = 46000
$ws.Range("I63").Value = 64000
$ws.Range("J63").Value = 10000
$ws.Range("K63").Value = 64000
$ws.Range("L63").Value = 10000
$ws.Range("M63").Value = -63314
$ws.Range("N63").Value = -11372

# row 66 (item 12528)
$ws.Range("H66").Value = 46000
$ws.Range("I66").Value = 64000
$ws.Range("J66").Value = 10000
$ws.Range("K66").Value = 320000
$ws.Range("L66").Value = 50000
$ws.Range("M66").Value = -316568
$ws.Range("N66").Value = -56864

# row 74 (item 44000)
$ws.Range("H74").Value = 5638.2573
$ws.Range("I74").Value = 2013.8
$ws.Range("J74").Value = 14699.4
$ws.Range("K74").Value = 2013.8
$ws.Range("L74").Value = 14699.4
$ws.Range("M74").Value = -1139.8
$ws.Range("N74").Value = -16447.4

# row 77 (item 44000)
$ws.Range("H77").Value = 5638.2573
$ws.Range("I77").Value = 2013.8
$ws.Range("J77").Value = 14699.4
$ws.Range("K77").Value = 10069
$ws.Range("L77").Value = 73497
$ws.Range("M77").Value = -5701
$ws.Range("N77").Value = -82233

$ws = $wb.Worksheets.Item("BSM")
# row 20 (item 14149)
$ws.Range("H20").Value = 2097.7222
$ws.Range("I20").Value = 1846.6666
$ws.Range("J20").Value = 2599.8333
$ws.Range("K20").Value = 1846.6666
$ws.Range("L20").Value = 2599.8333
$ws.Range("M20").Value = -1599.6666
$ws.Range("N20").Value = -3093.8333

# row 107 (item 27706)
$ws.Range("H107").Value = 1124.8235
$ws.Range("I107").Value = 690.5454999999999
$ws.Range("J107").Value = 1921
$ws.Range("K107").Value = 690.5454999999999
$ws.Range("L107").Value = 1921
$ws.Range("M107").Value = 1229.4545
$ws.Range("N107").Value = -5761

# row 134 (item 43998)
$ws.Range("H134").Value = 19232950
$ws.Range("I134").Value = 24392050
$ws.Range("J134").Value = 3579.6365
$ws.Range("K134").Value = 73176150
$ws.Range("L134").Value = 10738.9095
$ws.Range("M134").Value = -73173615
$ws.Range("N134").Value = -15808.9095

$ws = $wb.Worksheets.Item("CRP")
# row 31 (item 44023)
$ws.Range("H31").Value = 1753.3654
$ws.Range("I31").Value = 973.5897
$ws.Range("J31").Value = 4092.6924
$ws.Range("K31").Value = 973.5897
$ws.Range("L31").Value = 4092.6924
$ws.Range("M31").Value = -678.5897
$ws.Range("N31").Value = -4682.6924

# row 34 (item 44023)
$ws.Range("H34").Value = 1753.3654
$ws.Range("I34").Value = 973.5897
$ws.Range("J34").Value = 4092.6924
$ws.Range("K34").Value = 973.5897
$ws.Range("L34").Value = 4092.6924
$ws.Range("M34").Value = -771.5897
$ws.Range("N34").Value = -4496.6924

$ws = $wb.Worksheets.Item("CUL")
# row 5 (item 43974)
$ws.Range("H5").Value = 1076.7693
$ws.Range("I5").Value = 652.2857
$ws.Range("K5").Value = 1956.8571
$ws.Range("M5").Value = -1844.8571

# row 7 (item 4728)
$ws.Range("H7").Value = 350.2857
$ws.Range("I7").Value = 578
$ws.Range("J7").Value = 46.666668
$ws.Range("K7").Value = 1734
$ws.Range("L7").Value = 140.000004
$ws.Range("M7").Value = -1622
$ws.Range("N7").Value = -364.000004

# row 80 (item 12890)
$ws.Range("H80").Value = 1198.5714
$ws.Range("J80").Value = 1198.5714
$ws.Range("L80").Value = 3595.7142
$ws.Range("N80").Value = -5467.7142

# row 83 (item 12890)
$ws.Range("H83").Value = 1198.5714
$ws.Range("J83").Value = 1198.5714
$ws.Range("L83").Value = 10787.1426
$ws.Range("N83").Value = -20147.1426

# row 92 (item 19841)
$ws.Range("H92").Value = 823.125
$ws.Range("I92").Value = 800.5
$ws.Range("J92").Value = 845.75
$ws.Range("K92").Value = 2401.5
$ws.Range("L92").Value = 2537.25
$ws.Range("M92").Value = -1153.5
$ws.Range("N92").Value = -5033.25

# row 113 (item 27843)
$ws.Range("H113").Value = 12821070
$ws.Range("I113").Value = 569.9524
$ws.Range("J113").Value = 27778322
$ws.Range("K113").Value = 1709.8572
$ws.Range("L113").Value = 83334966
$ws.Range("M113").Value = 460.1428000000001
$ws.Range("N113").Value = -83339306

# row 122 (item 36078)
$ws.Range("H122").Value = 10835.2
$ws.Range("J122").Value = 17849.666
$ws.Range("L122").Value = 160646.994
$ws.Range("N122").Value = -165546.994

# row 132 (item 43972)
$ws.Range("H132").Value = 2701.6
$ws.Range("I132").Value = 2600
$ws.Range("J132").Value = 2727
$ws.Range("K132").Value = 23400
$ws.Range("L132").Value = 24543
$ws.Range("M132").Value = -20870
$ws.Range("N132").Value = -29603

# row 135 (item 43974)
$ws.Range("H135").Value = 1076.7693
$ws.Range("I135").Value = 652.2857
$ws.Range("K135").Value = 5870.571300000001
$ws.Range("M135").Value = -3335.571300000001

$ws = $wb.Worksheets.Item("GSM")
# row 113 (item 27710)
$ws.Range("H113").Value = 1922.5186
$ws.Range("I113").Value = 1398.8667
$ws.Range("J113").Value = 2577.0833
$ws.Range("K113").Value = 1398.8667
$ws.Range("L113").Value = 2577.0833
$ws.Range("M113").Value = 771.1333
$ws.Range("N113").Value = -6917.0833

$ws = $wb.Worksheets.Item("LTW")
# row 22 (item 5277)
$ws.Range("H22").Value = 7621.1763
$ws.Range("I22").Value = 1200
$ws.Range("J22").Value = 8997.143
$ws.Range("K22").Value = 1200
$ws.Range("L22").Value = 8997.143
$ws.Range("M22").Value = -905
$ws.Range("N22").Value = -9587.143

# row 27 (item 5277)
$ws.Range("H27").Value = 7621.1763
$ws.Range("I27").Value = 1200
$ws.Range("J27").Value = 8997.143
$ws.Range("K27").Value = 1200
$ws.Range("L27").Value = 8997.143
$ws.Range("M27").Value = -1093
$ws.Range("N27").Value = -9211.143

# row 82 (item 12565)
$ws.Range("H82").Value = 1990
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 1990
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 1990
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -2712

# row 85 (item 12565)
$ws.Range("H85").Value = 1990
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 1990
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 1990
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -4486

# row 121 (item 26320)
$ws.Range("H121").Value = 27622.857
$ws.Range("J121").Value = 27622.857
$ws.Range("L121").Value = 27622.857
$ws.Range("N121").Value = -31116.857

$ws = $wb.Worksheets.Item("WVR")
# row 113 (item 27752)
$ws.Range("H113").Value = 620.25
$ws.Range("I113").Value = 516.25
$ws.Range("J113").Value = 724.25
$ws.Range("K113").Value = 1548.75
$ws.Range("L113").Value = 2172.75
$ws.Range("M113").Value = 621.25
$ws.Range("N113").Value = -6512.75
